# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.456.08'
$ws.Cells.Item(2, 5).Value = '  +3.25%  '

$ws.Cells.Item(3, 4).Value = '2.625.66'
$ws.Cells.Item(3, 5).Value = '  +3.84%  '

$ws.Cells.Item(4, 4).Value = "'1.00"

$ws.Cells.Item(5, 4).Value = "'605.60"
$ws.Cells.Item(5, 5).Value = '  +1.64%  '

$ws.Cells.Item(6, 4).Value = "'179.55"
$ws.Cells.Item(6, 5).Value = '  +1.51%  '

$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$ws.Cells.Item(8, 5).Value = '  +1.02%  '

$ws.Cells.Item(9, 4).Value = '2.623.44'
$ws.Cells.Item(9, 5).Value = '  +3.73%  '

$ws.Cells.Item(10, 4).Value = "'0.167"
$ws.Cells.Item(10, 5).Value = '  +12.99%  '

$ws.Cells.Item(11, 5).Value = '  +0.25%  '

$ws.Cells.Item(12, 5).Value = '  +2.42%  '

$ws.Cells.Item(13, 4).Value = "'5.07"
$ws.Cells.Item(13, 5).Value = '  +1.51%  '

$ws.Cells.Item(14, 4).Value = '3.126.08'
$ws.Cells.Item(14, 5).Value = '  +5.92%  '

$ws.Cells.Item(15, 5).Value = '  +7.34%  '

$ws.Cells.Item(16, 4).Value = "'26.72"
$ws.Cells.Item(16, 5).Value = '  +1.87%  '

$ws.Cells.Item(17, 4).Value = '71.350.49'
$ws.Cells.Item(17, 5).Value = '  +3.48%  '

$ws.Cells.Item(18, 4).Value = '2.621.51'
$ws.Cells.Item(18, 5).Value = '  +3.60%  '

$ws.Cells.Item(19, 4).Value = "'381.20"
$ws.Cells.Item(19, 5).Value = '  +5.22%  '

$ws.Cells.Item(20, 4).Value = "'7.89"
$ws.Cells.Item(20, 5).Value = '  +4.61%  '

$ws.Cells.Item(21, 4).Value = '11.48'
$ws.Cells.Item(21, 5).Value = '  +3.38%  '

$ws.Cells.Item(22, 4).Value = "'4.13"
$ws.Cells.Item(22, 5).Value = '  +0.53%  '

$ws.Cells.Item(23, 4).Value = '1.99'
$ws.Cells.Item(23, 5).Value = '  +16.59%  '

$ws.Cells.Item(24, 4).Value = "'72.57"
$ws.Cells.Item(24, 5).Value = '  +2.49%  '

$ws.Cells.Item(25, 4).Value = "'4.43"
$ws.Cells.Item(25, 5).Value = '  +4.81%  '

$ws.Cells.Item(26, 5).Value = '  -0.01%  '

$ws.Cells.Item(27, 4).Value = '10.01'
$ws.Cells.Item(27, 5).Value = '  +10.70%  '

$ws.Cells.Item(28, 4).Value = '2.762.04'
$ws.Cells.Item(28, 5).Value = '  +4.04%  '

$ws.Cells.Item(29, 5).Value = '  +0.32%  '

$ws.Cells.Item(30, 4).Value = "'548.96"
$ws.Cells.Item(30, 5).Value = '  +5.13%  '

$ws.Cells.Item(31, 4).Value = '0.0₃0963'
$ws.Cells.Item(31, 5).Value = '  +7.25%  '

$ws.Cells.Item(32, 4).Value = "'8.09"
$ws.Cells.Item(32, 5).Value = '  +3.74%  '

$ws.Cells.Item(33, 5).Value = '  +7.03%  '

$ws.Cells.Item(34, 5).Value = '  +2.70%  '

$ws.Cells.Item(35, 5).Value = '  -0.20%  '

$ws.Cells.Item(36, 4).Value = "'166.49"

$ws.Cells.Item(37, 2).Value = 'EthereumClassic'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(37, 4).Value = '19.23'
$ws.Cells.Item(37, 5).Value = '  +3.70%  '

$ws.Cells.Item(38, 2).Value = 'Kaspa'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).Value = '0.115'
$ws.Cells.Item(38, 5).Value = '  -4.70%  '

$ws.Cells.Item(39, 4).Value = "'19.19"
$ws.Cells.Item(39, 5).Value = '  +2.67%  '

$ws.Cells.Item(40, 5).Value = '  +6.43%  '

$ws.Cells.Item(41, 4).Value = "'1.87"
$ws.Cells.Item(41, 5).Value = '  +5.33%  '

$ws.Cells.Item(42, 5).Value = '  +0.00%  '

$ws.Cells.Item(43, 4).Value = "'2.62"
$ws.Cells.Item(43, 5).Value = '  +9.01%  '

$ws.Cells.Item(44, 5).Value = '  +4.20%  '

$ws.Cells.Item(45, 5).Value = '  +2.01%  '

$ws.Cells.Item(46, 4).Value = "'39.94"
$ws.Cells.Item(46, 5).Value = '  +2.37%  '

$ws.Cells.Item(47, 4).Value = "'153.04"
$ws.Cells.Item(47, 5).Value = '  +1.03%  '

$ws.Cells.Item(48, 4).Value = "'3.64"
$ws.Cells.Item(48, 5).Value = '  +1.30%  '

$ws.Cells.Item(49, 5).Value = '  +3.40%  '

$ws.Cells.Item(50, 5).Value = '  +5.71%  '

$ws.Cells.Item(51, 5).Value = '  +3.32%  '
